# ChallengeCardData.xlsx edit
#
# The Relic/"遗物" challenge card (row 4, column C) gets its effect text
# extended with a new win-condition line ("胜利条件：《收藏家》牌进入墓地。").
# In the canonical OOXML this shows up as the old shared-string entry being
# removed and a new (longer) one appended at the end of the shared-strings
# table, which shuffles every other row's <v> string index even though their
# actual text is unchanged -- so here we just rewrite the one cell whose
# wording really changed; Excel takes care of the shared-string bookkeeping.
#
# Row 4 also grows from two wrapped lines to three, so its row height grows
# from 57 to 71.25 (matching the other three-line rows in the sheet).
#
# Finally the sheet's active view is updated: the frozen "topLeftCell"
# scroll position is cleared and the selected cell moves from G5 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRelicEffect = @"
DL2：房间尺寸加1，手牌基数加1。遗物牌使用后横置，且不会在回合结束时被复位。<br>
DL3：房间尺寸加1，手牌基数加1。首次重整后，从额外牌堆将1张《收藏家》牌洗入主牌堆。<br>
胜利条件：《收藏家》牌进入墓地。
"@

$ws.Range("C4").Value = $newRelicEffect
$ws.Rows(4).RowHeight = 71.25

$ws.Range("C5").Select()
